# Auto-generated from unified OOXML diff.
# Updates coin Price (D) and Volume(1h) (E) values, and reorders some
# Coin/Link (B/C) rows for the "Updated symbol list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.49"
$ws.Range("E2").Value = "'1.58%"
$ws.Range("D3").Value = "'45.85"
$ws.Range("E3").Value = "'4.25%"
$ws.Range("D4").Value = "'5.652"
$ws.Range("E4").Value = "'2.85%"
$ws.Range("D5").Value = "'0.08390"
$ws.Range("E5").Value = "'4.67%"
$ws.Range("D6").Value = "'2.032"
$ws.Range("E6").Value = "'0.84%"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = "'4.483"
$ws.Range("E7").Value = "'3.83%"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9944"
$ws.Range("E8").Value = "'4.84%"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = "'2.578"
$ws.Range("E9").Value = "'0.36%"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.1150"
$ws.Range("E10").Value = "'2.24%"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1930"
$ws.Range("E11").Value = "'3.79%"
$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").Value = "'10.38"
$ws.Range("E12").Value = "'-2.25%"
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = "'0.09958"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = "'0.04673"
$ws.Range("E14").Value = "'1.78%"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.1059"
$ws.Range("E15").Value = "'-0.72%"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001276"
$ws.Range("E16").Value = "'-0.05%"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = "'0.006038"
$ws.Range("E17").Value = "'1.94%"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = "'3.376"
$ws.Range("E18").Value = "'0.68%"
$ws.Range("E19").Value = "'-3.14%"
$ws.Range("D20").Value = "'0.1402"
$ws.Range("E20").Value = "'-0.45%"
$ws.Range("D21").Value = "'0.2654"
$ws.Range("E21").Value = "'4.36%"
$ws.Range("D22").Value = "'0.04224"
$ws.Range("E22").Value = "'3.72%"
$ws.Range("D23").Value = "'0.001310"
$ws.Range("E23").Value = "'3.95%"
$ws.Range("D24").Value = "'0.004640"
$ws.Range("E24").Value = "'7.37%"
$ws.Range("E25").Value = "'10.79%"
$ws.Range("D26").Value = "'0.0003749"
$ws.Range("E26").Value = "'0.26%"
$ws.Range("D38").Value = "'0.02787"
$ws.Range("E38").Value = "'8.35%"
$ws.Range("D39").Value = "'0.05758"
$ws.Range("E39").Value = "'1.37%"
$ws.Range("D40").Value = "'0.007744"
$ws.Range("E40").Value = "'2.70%"
$ws.Range("D41").Value = "'0.1434"
$ws.Range("E41").Value = "'2.68%"
$ws.Range("D42").Value = "'0.007247"
$ws.Range("E42").Value = "'-4.49%"
$ws.Range("D43").Value = "'0.002119"
$ws.Range("E43").Value = "'5.29%"
$ws.Range("D44").Value = "'0.009053"
$ws.Range("E44").Value = "'8.04%"
$ws.Range("D45").Value = "'0.3409"
$ws.Range("D46").Value = "'0.00007355"
$ws.Range("E46").Value = "'3.69%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.37%"
$ws.Range("D48").Value = "'0.0005815"
$ws.Range("E48").Value = "'0.06%"
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").Value = "'0.003507"
$ws.Range("E49").Value = "'-0.64%"
$ws.Range("B50").Value = 'BOLO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D50").Value = "'0.003509"
$ws.Range("E50").Value = "'4.01%"
$ws.Range("D51").Value = "'0.00002105"
$ws.Range("E51").Value = "'0.37%"
